$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 8

$ws.Range("B8").Select()
